# Insert two new data rows right after the header block / before the current row 324
# (pushing the existing rows 324..403 down to 326..405), then fill in the two new rows
# with the new weekly price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two whole rows at 324-325; everything currently at row 324 onward shifts down by 2.
$ws.Range("A324:A325").EntireRow.Insert()

# Common (constant) columns for this product/market sheet.
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonE = 13
$commonF = "Fruta"
$commonG = 100101
$commonH = "Berries"
$commonI = 100101001
$commonJ = "Arándano (blue)"
$commonK = "Sin especificar"
$commonQ = "`$/bandeja 2 kilos"
$commonT = 2

# ---- New row 324 ----
$ws.Cells.Item(324, 1).Value = $commonA
$ws.Cells.Item(324, 2).Value = $commonB
$ws.Cells.Item(324, 3).Value = $commonC
$ws.Cells.Item(324, 4).Value = 45275
$ws.Cells.Item(324, 5).Value = $commonE
$ws.Cells.Item(324, 6).Value = $commonF
$ws.Cells.Item(324, 7).Value = $commonG
$ws.Cells.Item(324, 8).Value = $commonH
$ws.Cells.Item(324, 9).Value = $commonI
$ws.Cells.Item(324, 10).Value = $commonJ
$ws.Cells.Item(324, 11).Value = $commonK
$ws.Cells.Item(324, 12).Value = "Especial"
$ws.Cells.Item(324, 13).Value = 440
$ws.Cells.Item(324, 14).Value = 3800
$ws.Cells.Item(324, 15).Value = 3800
$ws.Cells.Item(324, 16).Value = 3800
$ws.Cells.Item(324, 17).Value = $commonQ
$ws.Cells.Item(324, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(324, 19).Value = 1900
$ws.Cells.Item(324, 20).Value = $commonT

# ---- New row 325 ----
$ws.Cells.Item(325, 1).Value = $commonA
$ws.Cells.Item(325, 2).Value = $commonB
$ws.Cells.Item(325, 3).Value = $commonC
$ws.Cells.Item(325, 4).Value = 45275
$ws.Cells.Item(325, 5).Value = $commonE
$ws.Cells.Item(325, 6).Value = $commonF
$ws.Cells.Item(325, 7).Value = $commonG
$ws.Cells.Item(325, 8).Value = $commonH
$ws.Cells.Item(325, 9).Value = $commonI
$ws.Cells.Item(325, 10).Value = $commonJ
$ws.Cells.Item(325, 11).Value = $commonK
$ws.Cells.Item(325, 12).Value = "Primera"
$ws.Cells.Item(325, 13).Value = 520
$ws.Cells.Item(325, 14).Value = 3400
$ws.Cells.Item(325, 15).Value = 3400
$ws.Cells.Item(325, 16).Value = 3400
$ws.Cells.Item(325, 17).Value = $commonQ
$ws.Cells.Item(325, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(325, 19).Value = 1700
$ws.Cells.Item(325, 20).Value = $commonT
